$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 2583
$ws.Range("J3").Value = 2638
$ws.Range("I4").Value = 1757
$ws.Range("J4").Value = 606
$ws.Range("J5").Value = 205
$ws.Range("H6").Value = 7917
$ws.Range("J6").Value = 3277
$ws.Range("H7").Value = 26006
$ws.Range("I7").Value = 26204
$ws.Range("J7").Value = 9309

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 96
$ws.Range("J3").Value = 111
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 315

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 43
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 71
$ws.Range("J3").Value = 135
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 66
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 241

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 72
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 279
$ws.Range("J8").Value = 581
$ws.Range("J11").Value = 132
$ws.Range("H12").Value = 49
$ws.Range("J19").Value = 300
$ws.Range("J23").Value = 97
$ws.Range("J25").Value = 53
$ws.Range("J29").Value = 540
$ws.Range("J31").Value = 70
$ws.Range("J33").Value = 378
$ws.Range("J34").Value = 48
$ws.Range("J36").Value = 139
$ws.Range("J37").Value = 315
$ws.Range("J41").Value = 62
$ws.Range("J42").Value = 361
$ws.Range("J43").Value = 85
$ws.Range("J46").Value = 30
$ws.Range("J47").Value = 79
$ws.Range("J48").Value = 92
$ws.Range("J49").Value = 59
$ws.Range("J50").Value = 52
$ws.Range("J51").Value = 125
$ws.Range("J54").Value = 186
$ws.Range("I63").Value = 208
$ws.Range("J63").Value = 48
$ws.Range("J65").Value = 241
$ws.Range("J67").Value = 338
$ws.Range("J72").Value = 35
$ws.Range("J73").Value = 82
$ws.Range("J74").Value = 14
$ws.Range("J76").Value = 132
$ws.Range("J77").Value = 75
$ws.Range("J78").Value = 129
$ws.Range("J85").Value = 433
$ws.Range("J86").Value = 56
$ws.Range("J90").Value = 103
$ws.Range("J91").Value = 105
$ws.Range("J95").Value = 143
$ws.Range("J96").Value = 108
$ws.Range("J98").Value = 54
$ws.Range("J99").Value = 128
$ws.Range("H101").Value = 26006
$ws.Range("I101").Value = 26204
$ws.Range("J101").Value = 9309

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 53
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 113
$ws.Range("J5").Value = 16
$ws.Range("J7").Value = 378

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J2").Value = 13
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J4").Value = 15
$ws.Range("J7").Value = 186

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 159
$ws.Range("J3").Value = 178
$ws.Range("J5").Value = 24
$ws.Range("J7").Value = 540

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 83
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 104
$ws.Range("J3").Value = 164
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 433

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 90

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J2").Value = 17
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 72
$ws.Range("J3").Value = 77
$ws.Range("J7").Value = 361

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 105

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 139

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 46
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 30
$ws.Range("J5").Value = 1
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 183
$ws.Range("J3").Value = 189
$ws.Range("J7").Value = 581

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 27
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 33
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 36
$ws.Range("J5").Value = 4
$ws.Range("J7").Value = 125

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 24
$ws.Range("J3").Value = 24
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 92
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 279

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("H6").Value = 20
$ws.Range("H7").Value = 49

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 14
